$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1814584051120676
    "C2" = 0.03523900328895735
    "D2" = 0.02953349150218543
    "F2" = 0.6783337081030183
    "G2" = 0.002424109022594731
    "I2" = 0.5378927371314219
    "K2" = 0.1867389260986982
    "M2" = 0.8848238073444321
    "O2" = 2.319661342300435
    "B3" = 0.1587103975088695
    "C3" = 0.03244207139204747
    "D3" = 0.02733554122478665
    "F3" = 0.6766109037018921
    "G3" = 0.002426378543332483
    "I3" = 0.541628343323282
    "K3" = 0.1630004039512585
    "M3" = 0.7879428552989793
    "O3" = 2.326631108638026
    "B4" = 0.1447039711759714
    "C4" = 0.03070913901017747
    "D4" = 0.02597297157370093
    "F4" = 0.6759675610600695
    "G4" = 0.002427845208358282
    "I4" = 0.5442145161492213
    "K4" = 0.1483606356874816
    "M4" = 0.7289846298061065
    "O4" = 2.33224154640348
    "B5" = 0.1389869409116784
    "C5" = 0.02999909472234208
    "D5" = 0.02541448169869653
    "F5" = 0.6758095288863473
    "G5" = 0.002428461345263518
    "I5" = 0.5453418639107781
    "K5" = 0.1423792007681044
    "M5" = 0.7050851701381333
    "O5" = 2.334862161224081
    "B6" = 0.138037086257782
    "C6" = 0.02988096113245575
    "D6" = 0.02532155093458499
    "F6" = 0.6757895747821294
    "G6" = 0.00242856477103823
    "I6" = 0.5455334934349771
    "K6" = 0.1413850610758658
    "M6" = 0.7011241485861888
    "O6" = 2.335317493456614
    "B7" = 0.1446269062582672
    "C7" = 0.03069957865056239
    "D7" = 0.02596545261342698
    "F7" = 0.6759650082606328
    "G7" = 0.002427853442917393
    "I7" = 0.5442294226351692
    "K7" = 0.1482800304571015
    "M7" = 0.7286618090380301
    "O7" = 2.332275535801301
    "B8" = 0.1736233051890679
    "C8" = 0.03427789736544185
    "D8" = 0.02877836505206233
    "F8" = 0.6776536084636007
    "G8" = 0.002424876403141429
    "I8" = 0.539120012693914
    "K8" = 0.1785675277078411
    "M8" = 0.8513061669071931
    "O8" = 2.321788047846042
    "B9" = 0.230154936194765
    "C9" = 0.0411686501438453
    "D9" = 0.03418955646281319
    "F9" = 0.6842586111419919
    "G9" = 0.002419616271946647
    "I9" = 0.5314255314897345
    "K9" = 0.2374307614200291
    "M9" = 1.096281449235775
    "O9" = 2.311801655360995
    "B10" = 0.2714640548322507
    "C10" = 0.04615135105071033
    "D10" = 0.03809935955373334
    "F10" = 0.6911279447262686
    "G10" = 0.002416100069946976
    "I10" = 0.5271959987068833
    "K10" = 0.2803301189262015
    "M10" = 1.279411530931768
    "O10" = 2.310943630725461
    "B11" = 0.2902031564296692
    "C11" = 0.0484001577421509
    "D11" = 0.03986337233384774
    "F11" = 0.6946928785154469
    "G11" = 0.002414575291030462
    "I11" = 0.5255823686016505
    "K11" = 0.2997658457399552
    "M11" = 1.363500004184417
    "O11" = 2.311966755193481
    "B12" = 0.2972911194000005
    "C12" = 0.04924909567914426
    "D12" = 0.0405292229713865
    "F12" = 0.6961062341132944
    "G12" = 0.002414008584683595
    "I12" = 0.5250160755849933
    "K12" = 0.3071137358050748
    "M12" = 1.39546210750413
    "O12" = 2.312557928237112
    "B13" = 0.2957649701682783
    "C13" = 0.04906638006583819
    "D13" = 0.0403859161630109
    "F13" = 0.6957990214503695
    "G13" = 0.002414130160305283
    "I13" = 0.5251360445999538
    "K13" = 0.3055317790277456
    "M13" = 1.388573059504324
    "O13" = 2.312421539449815
    "B14" = 0.2907864535159206
    "C14" = 0.0484700536134568
    "D14" = 0.0399181954769503
    "F14" = 0.694807884986794
    "G14" = 0.002414528453938176
    "I14" = 0.5255348817262195
    "K14" = 0.3003706048469326
    "M14" = 1.366127106222493
    "O14" = 2.312011305342963
    "B15" = 0.2877358931752951
    "C15" = 0.04810444104242606
    "D15" = 0.03963142262104924
    "F15" = 0.6942090443011253
    "G15" = 0.002414773811249815
    "I15" = 0.5257850127432064
    "K15" = 0.2972076569866715
    "M15" = 1.352394107609967
    "O15" = 2.311786572312002
    "B16" = 0.2702383035658613
    "C16" = 0.04600402111952917
    "D16" = 0.03798377986873902
    "F16" = 0.6909038324870096
    "G16" = 0.002416201217292174
    "I16" = 0.5273077094255889
    "K16" = 0.2790583010580292
    "M16" = 1.273932493797773
    "O16" = 2.310905250741229
    "B17" = 0.2594902212931629
    "C17" = 0.04471086038763872
    "D17" = 0.03696923787300221
    "F17" = 0.6889889786439554
    "G17" = 0.002417095992369112
    "I17" = 0.5283214203003723
    "K17" = 0.2679035085669739
    "M17" = 1.226004121372725
    "O17" = 2.310726950357065
    "B18" = 0.2533032888746902
    "C18" = 0.04396539267385435
    "D18" = 0.03638433095331806
    "F18" = 0.687929017031756
    "G18" = 0.002417617683426894
    "I18" = 0.5289336897387535
    "K18" = 0.2614801274172862
    "M18" = 1.198510122043146
    "O18" = 2.310757407478775
    "B19" = 0.2512076720070695
    "C19" = 0.04371270440803698
    "D19" = 0.03618605778225259
    "F19" = 0.687577241227558
    "G19" = 0.002417795529799686
    "I19" = 0.5291460068775606
    "K19" = 0.259304020858167
    "M19" = 1.189213469941166
    "O19" = 2.310790549689841
    "B20" = 0.2606348862450432
    "C20" = 0.04484869350622489
    "D20" = 0.03707737964514735
    "F20" = 0.6891885313212214
    "G20" = 0.002417000013921165
    "I20" = 0.5282104850996312
    "K20" = 0.2690917294354449
    "M20" = 1.231098557753398
    "O20" = 2.310732161118807
    "B21" = 0.2922489895942704
    "C21" = 0.04864528125224865
    "D21" = 0.04005563490546393
    "F21" = 0.6950972844128174
    "G21" = 0.002414411175651932
    "I21" = 0.5254165180487931
    "K21" = 0.3018868970416122
    "M21" = 1.372716721815863
    "O21" = 2.312126267755872
    "B22" = 0.3128629976328625
    "C22" = 0.05111117180497615
    "D22" = 0.0419895858111019
    "F22" = 0.6993285504660349
    "G22" = 0.002412781528151318
    "I22" = 0.5238514131251719
    "K22" = 0.3232502342016517
    "M22" = 1.465973117963543
    "O22" = 2.314225166671719
    "B23" = 0.3018654572057926
    "C23" = 0.04979651164420318
    "D23" = 0.04095856000321874
    "F23" = 0.697036391563401
    "G23" = 0.002413645618973241
    "I23" = 0.5246628262603394
    "K23" = 0.3118548335542641
    "M23" = 1.416133873431065
    "O23" = 2.312996102028961
    "B24" = 0.2601174072129595
    "C24" = 0.04478638541584701
    "D24" = 0.03702849383352458
    "F24" = 0.6890981861284899
    "G24" = 0.002417043383164786
    "I24" = 0.5282605471092623
    "K24" = 0.2685545669507121
    "M24" = 1.228795174665777
    "O24" = 2.310729391215688
    "B25" = 0.2148994068256513
    "C25" = 0.03931834530266087
    "D25" = 0.03273710832489485
    "F25" = 0.6821182729056474
    "G25" = 0.002420977814320246
    "I25" = 0.5332575429103308
    "K25" = 0.2215659822369389
    "M25" = 1.029487586080663
    "O25" = 2.313367598994773
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
